$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: seed rows 722-771 by duplicating the last existing row (721) so
#     formatting/styles (s="1" / s="2" date col) carry over exactly. ---
$ws.Range("A721:H721").Copy($ws.Range("A722:H771"))

# --- Step 2: write the real values for columns B..H, row by row (722->771),
#     left to right -- this is the order new text actually entered the sheet. ---
# Row 722
$ws.Range("B722").Value = 'SAP'
$ws.Range("C722").Value = 'XS-PTS-0907'
$ws.Range("D722").Value = '70900.067 SA/MP Test Socket Pogo Pin X1893'
$ws.Range("E722").Value = '400pcs'
$ws.Range("F722").Value = 'Fishes'
$ws.Range("G722").Value = '13/10/2025'
$ws.Range("H722").Value = 1776

# Row 723
$ws.Range("B723").Value = 'SAP'
$ws.Range("C723").Value = 11155497
$ws.Range("D723").Value = 'PTS-1138 6K-9090-HS01  LF Vacuum Adaptor'
$ws.Range("E723").Value = '100pcs'
$ws.Range("F723").Value = 'Fishes'
$ws.Range("G723").Value = '13/10/2025'
$ws.Range("H723").Value = 3510

# Row 724
$ws.Range("B724").Value = 'SAP'
$ws.Range("C724").Value = 11155143
$ws.Range("D724").Value = 'PTS-1136 300-001519-015 Semiconductor Primary Pogo pin X2637'
$ws.Range("E724").Value = '1400pcs'
$ws.Range("F724").Value = 'Fishes'
$ws.Range("G724").Value = '13/10/2025'
$ws.Range("H724").Value = 2604

# Row 725
$ws.Range("B725").Value = 'SAP'
$ws.Range("C725").Value = 11155143
$ws.Range("D725").Value = 'PTS-1136 300-001519-015 Semiconductor Primary Pogo pin X2637'
$ws.Range("E725").Value = '200pcs'
$ws.Range("F725").Value = 'Fishes'
$ws.Range("G725").Value = '13/10/2025'
$ws.Range("H725").Value = 372

# Row 726
$ws.Range("B726").Value = 'SAP'
$ws.Range("C726").Value = 'XS-PTS-1043'
$ws.Range("D726").Value = '10416.113 SAMTEC CABLE 1.2M X1767'
$ws.Range("E726").Value = '3pcs'
$ws.Range("F726").Value = 'Fishes'
$ws.Range("G726").Value = '13/10/2025'
$ws.Range("H726").Value = 2098.7600000000002

# Row 727
$ws.Range("B727").Value = 'SAP'
$ws.Range("C727").Value = 11155949
$ws.Range("D727").Value = 'PTS-1146 70902.631 X2544 ALIGNER & NEST'
$ws.Range("E727").Value = '4pcs'
$ws.Range("F727").Value = 'Fishes'
$ws.Range("G727").Value = '13/10/2025'
$ws.Range("H727").Value = 11557.92

# Row 728
$ws.Range("B728").Value = 'SAP'
$ws.Range("C728").Value = 11156537
$ws.Range("D728").Value = 'PTS-1150 300-01698-010 Semiconductor Secondary Pogo Pin X2544'
$ws.Range("E728").Value = '400pcs'
$ws.Range("F728").Value = 'Fishes'
$ws.Range("G728").Value = '13/10/2025'
$ws.Range("H728").Value = 936

# Row 729
$ws.Range("B729").Value = 'SAP'
$ws.Range("C729").Value = 'XS-PTS-0867'
$ws.Range("D729").Value = '6K-57084-H062 Hyperspace Semicon Socket'
$ws.Range("E729").Value = '1pcs'
$ws.Range("F729").Value = 'Sihl'
$ws.Range("G729").Value = '14/10/2025'
$ws.Range("H729").Value = 190.99

# Row 730
$ws.Range("B730").Value = 'SAP'
$ws.Range("C730").Value = 'XS-PTS-0502'
$ws.Range("D730").Value = 'Hyperspace Socket Model : 6K-76235-H04'
$ws.Range("E730").Value = '1pcs'
$ws.Range("F730").Value = 'Sihl'
$ws.Range("G730").Value = '14/10/2025'
$ws.Range("H730").Value = 191.35

# Row 731
$ws.Range("B731").Value = 'SAP'
$ws.Range("C731").Value = 'XS-PTS-0994'
$ws.Range("D731").Value = 'P/N: PJWBB790-28-00S ISC P-PIN'
$ws.Range("E731").Value = '400pcs'
$ws.Range("F731").Value = 'Sihl'
$ws.Range("G731").Value = '14/10/2025'
$ws.Range("H731").Value = 720

# Row 732
$ws.Range("B732").Value = 'SAP'
$ws.Range("C732").Value = 'XS-SPE-0096'
$ws.Range("D732").Value = 'Wired motor for RRU300 (mod. L20A18S0604'
$ws.Range("E732").Value = '1pcs'
$ws.Range("F732").Value = 'Lisa'
$ws.Range("G732").Value = '14/10/2025'
$ws.Range("H732").Value = 331.6

# Row 733
$ws.Range("B733").Value = 'SAP'
$ws.Range("C733").Value = 11155797
$ws.Range("D733").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E733").Value = '100pcs'
$ws.Range("F733").Value = 'Fishes'
$ws.Range("G733").Value = '14/10/2025'
$ws.Range("H733").Value = 2016

# Row 734
$ws.Range("B734").Value = 'SAP'
$ws.Range("C734").Value = 11155497
$ws.Range("D734").Value = 'PTS-1138 6K-9090-HS01  LF Vacuum Adaptor'
$ws.Range("E734").Value = '100pcs'
$ws.Range("F734").Value = 'Fishes'
$ws.Range("G734").Value = '14/10/2025'
$ws.Range("H734").Value = 3510

# Row 735
$ws.Range("B735").Value = 'SAP'
$ws.Range("C735").Value = 11151248
$ws.Range("D735").Value = 'PTS-1071 TW.50.15.FI.0S.151.00 X2637'
$ws.Range("E735").Value = '4pcs'
$ws.Range("F735").Value = 'Fishes'
$ws.Range("G735").Value = '14/10/2025'
$ws.Range("H735").Value = 2222.2399999999998

# Row 736
$ws.Range("B736").Value = 'SAP'
$ws.Range("C736").Value = 11151246
$ws.Range("D736").Value = 'PTS-1069 TW.50.15.FI.0S.150.00 X2637'
$ws.Range("E736").Value = '2pcs'
$ws.Range("F736").Value = 'Fishes'
$ws.Range("G736").Value = '14/10/2025'
$ws.Range("H736").Value = 1111.1199999999999

# Row 737
$ws.Range("B737").Value = 'SAP'
$ws.Range("C737").Value = 11155143
$ws.Range("D737").Value = 'PTS-1136 300-001519-015 Semiconductor Primary Pogo pin X2637'
$ws.Range("E737").Value = '500pcs'
$ws.Range("F737").Value = 'Fishes'
$ws.Range("G737").Value = '14/10/2025'
$ws.Range("H737").Value = 930

# Row 738
$ws.Range("B738").Value = 'SAP'
$ws.Range("C738").Value = 11151237
$ws.Range("D738").Value = 'PTS-1060 70192.692 PNP RUBBER TIP X2637'
$ws.Range("E738").Value = '200pcs'
$ws.Range("F738").Value = 'Fishes'
$ws.Range("G738").Value = '14/10/2025'
$ws.Range("H738").Value = 1354

# Row 739
$ws.Range("B739").Value = 'SAP'
$ws.Range("C739").Value = 'XS-PTS-0876'
$ws.Range("D739").Value = 'HX 2067 HPN CRC CO CONTACT CLEANER'
$ws.Range("E739").Value = '4pcs'
$ws.Range("F739").Value = 'Fishes'
$ws.Range("G739").Value = '15/10/2025'
$ws.Range("H739").Value = 133.63999999999999

# Row 740
$ws.Range("B740").Value = 'SAP'
$ws.Range("C740").Value = 'XS-PTS-1025'
$ws.Range("D740").Value = '10618.38 PRIMARY PCB BOARD ASSY X1767'
$ws.Range("E740").Value = '4pcs'
$ws.Range("F740").Value = 'Fishes'
$ws.Range("G740").Value = '15/10/2025'
$ws.Range("H740").Value = 421

# Row 741
$ws.Range("B741").Value = 'SAP'
$ws.Range("C741").Value = 'XS-PTS-1027'
$ws.Range("D741").Value = '10618.396 DUT DAUGHTER CARD X1767'
$ws.Range("E741").Value = '4pcs'
$ws.Range("F741").Value = 'Fishes'
$ws.Range("G741").Value = '15/10/2025'
$ws.Range("H741").Value = 7703.45

# Row 742
$ws.Range("B742").Value = 'SAP'
$ws.Range("C742").Value = 'XS-PTS-0907'
$ws.Range("D742").Value = '70900.067 SA/MP Test Socket Pogo Pin X1893'
$ws.Range("E742").Value = '400pcs'
$ws.Range("F742").Value = 'Fishes'
$ws.Range("G742").Value = '15/10/2025'
$ws.Range("H742").Value = 1776

# Row 743
$ws.Range("B743").Value = 'SAP'
$ws.Range("C743").Value = 'XS-PTS-0963'
$ws.Range("D743").Value = 'TW.50.1A.FI.0R.030.02 Aligner X1893'
$ws.Range("E743").Value = '2pcs'
$ws.Range("F743").Value = 'Fishes'
$ws.Range("G743").Value = '15/10/2025'
$ws.Range("H743").Value = 132.86000000000001

# Row 744
$ws.Range("B744").Value = 'SAP'
$ws.Range("C744").Value = 'XS-PTS-0971'
$ws.Range("D744").Value = 'TW.50.15.FI.0R.061.00 STT B ADAPTOR X189'
$ws.Range("E744").Value = '3pcs'
$ws.Range("F744").Value = 'Fishes'
$ws.Range("G744").Value = '15/10/2025'
$ws.Range("H744").Value = 1785.58

# Row 745
$ws.Range("B745").Value = 'SAP'
$ws.Range("C745").Value = 'XS-PTS-0972'
$ws.Range("D745").Value = '70192.496 RUBBER TIP DARK LEAKAGE X1893'
$ws.Range("E745").Value = '6pcs'
$ws.Range("F745").Value = 'Fishes'
$ws.Range("G745").Value = '15/10/2025'
$ws.Range("H745").Value = 1006.88

# Row 746
$ws.Range("B746").Value = 'SAP'
$ws.Range("C746").Value = 11151248
$ws.Range("D746").Value = 'PTS-1071 TW.50.15.FI.0S.151.00 X2637'
$ws.Range("E746").Value = '3pcs'
$ws.Range("F746").Value = 'Fishes'
$ws.Range("G746").Value = '15/10/2025'
$ws.Range("H746").Value = 1666.68

# Row 747
$ws.Range("B747").Value = 'SAP'
$ws.Range("C747").Value = 11151246
$ws.Range("D747").Value = 'PTS-1069 TW.50.15.FI.0S.150.00 X2637'
$ws.Range("E747").Value = '3pcs'
$ws.Range("F747").Value = 'Fishes'
$ws.Range("G747").Value = '15/10/2025'
$ws.Range("H747").Value = 1666.68

# Row 748
$ws.Range("B748").Value = 'SAP'
$ws.Range("C748").Value = 11151237
$ws.Range("D748").Value = 'PTS-1060 70192.692 PNP RUBBER TIP X2637'
$ws.Range("E748").Value = '260pcs'
$ws.Range("F748").Value = 'Fishes'
$ws.Range("G748").Value = '15/10/2025'
$ws.Range("H748").Value = 1760.2

# Row 749
$ws.Range("B749").Value = 'SAP'
$ws.Range("C749").Value = 11155143
$ws.Range("D749").Value = 'PTS-1136 300-001519-015 Semiconductor Primary Pogo pin X2637'
$ws.Range("E749").Value = '500pcs'
$ws.Range("F749").Value = 'Fishes'
$ws.Range("G749").Value = '15/10/2025'
$ws.Range("H749").Value = 930

# Row 750
$ws.Range("B750").Value = 'SAP'
$ws.Range("C750").Value = 'XS-PTS-0752'
$ws.Range("D750").Value = 'HX 1483 HPN OMRON EE-SX952-W-1M SENSOR X1893'
$ws.Range("E750").Value = '4pcs'
$ws.Range("F750").Value = 'Fishes'
$ws.Range("G750").Value = '15/10/2025'
$ws.Range("H750").Value = 60

# Row 751
$ws.Range("B751").Value = 'SAP'
$ws.Range("C751").Value = 11156537
$ws.Range("D751").Value = 'PTS-1150 300-01698-010 Semiconductor Secondary Pogo Pin X2544'
$ws.Range("E751").Value = '100pcs'
$ws.Range("F751").Value = 'Fishes'
$ws.Range("G751").Value = '15/10/2025'
$ws.Range("H751").Value = 234

# Row 752
$ws.Range("B752").Value = 'SAP'
$ws.Range("C752").Value = 11156536
$ws.Range("D752").Value = 'PTS-1149 40Ways Ribbon Flat Cable (1.2M) '
$ws.Range("E752").Value = '10pcs'
$ws.Range("F752").Value = 'Fishes'
$ws.Range("G752").Value = '15/10/2025'
$ws.Range("H752").Value = 194.66

# Row 753
$ws.Range("B753").Value = 'SAP'
$ws.Range("C753").Value = 'XS-MISC-0051'
$ws.Range("D753").Value = 'Heat Shrinkable Tube Inside Diameter 2mm'
$ws.Range("E753").Value = '2pcs'
$ws.Range("F753").Value = 'Fishes'
$ws.Range("G753").Value = '15/10/2025'
$ws.Range("H753").Value = 0.75

# Row 754
$ws.Range("B754").Value = 'SAP'
$ws.Range("C754").Value = 11156541
$ws.Range("D754").Value = 'PTS-1154 70902.708 TOP NEST X2544'
$ws.Range("E754").Value = '2pcs'
$ws.Range("F754").Value = 'Fishes'
$ws.Range("G754").Value = '15/10/2025'
$ws.Range("H754").Value = 784.32

# Row 755
$ws.Range("B755").Value = 'SAP'
$ws.Range("C755").Value = 11155947
$ws.Range("D755").Value = 'PTS-1144 800.403.00 X2544 MYCROFTL MP'
$ws.Range("E755").Value = '2pcs'
$ws.Range("F755").Value = 'Fishes'
$ws.Range("G755").Value = '15/10/2025'
$ws.Range("H755").Value = 752.32

# Row 756
$ws.Range("B756").Value = 'SAP'
$ws.Range("C756").Value = 11154690
$ws.Range("D756").Value = 'PTS-1133 TW.50.1B.98.JV.001.02 SHUTTLE X2544'
$ws.Range("E756").Value = '2pcs'
$ws.Range("F756").Value = 'Fishes'
$ws.Range("G756").Value = '15/10/2025'
$ws.Range("H756").Value = 833.34

# Row 757
$ws.Range("B757").Value = 'SAP'
$ws.Range("C757").Value = 'XS-PTS-1030'
$ws.Range("D757").Value = 'TW.50.1A.00.02.016.00 BRASS SLIDE X1767'
$ws.Range("E757").Value = '10pcs'
$ws.Range("F757").Value = 'Fishes'
$ws.Range("G757").Value = '15/10/2025'
$ws.Range("H757").Value = 135.30000000000001

# Row 758
$ws.Range("B758").Value = 'SAP'
$ws.Range("C758").Value = 'XS-PTS-1032'
$ws.Range("D758").Value = 'X1767 800.404.C3 2nd PG BLOCK 8 STAGE(S)'
$ws.Range("E758").Value = '4pcs'
$ws.Range("F758").Value = 'Fishes'
$ws.Range("G758").Value = '15/10/2025'
$ws.Range("H758").Value = 1850.99

# Row 759
$ws.Range("B759").Value = 'SAP'
$ws.Range("C759").Value = 11155797
$ws.Range("D759").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E759").Value = '200pcs'
$ws.Range("F759").Value = 'Fishes'
$ws.Range("G759").Value = '15/10/2025'
$ws.Range("H759").Value = 4032

# Row 760
$ws.Range("B760").Value = 'SAP'
$ws.Range("C760").Value = 11155497
$ws.Range("D760").Value = 'PTS-1138 6K-9090-HS01  LF Vacuum Adaptor'
$ws.Range("E760").Value = '145pcs'
$ws.Range("F760").Value = 'Fishes'
$ws.Range("G760").Value = '15/10/2025'
$ws.Range("H760").Value = 5089.5

# Row 761
$ws.Range("B761").Value = 'SAP'
$ws.Range("C761").Value = 11155143
$ws.Range("D761").Value = 'PTS-1136 300-001519-015 Semiconductor Primary Pogo pin X2637'
$ws.Range("E761").Value = '500pcs'
$ws.Range("F761").Value = 'Fishes'
$ws.Range("G761").Value = '15/10/2025'
$ws.Range("H761").Value = 930

# Row 762
$ws.Range("B762").Value = 'SAP'
$ws.Range("C762").Value = 11151237
$ws.Range("D762").Value = 'PTS-1060 70192.692 PNP RUBBER TIP X2637'
$ws.Range("E762").Value = '540pcs'
$ws.Range("F762").Value = 'Fishes'
$ws.Range("G762").Value = '15/10/2025'
$ws.Range("H762").Value = 3655.8

# Row 763
$ws.Range("B763").Value = 'SAP'
$ws.Range("C763").Value = 11156539
$ws.Range("D763").Value = 'PTS-1152 TW.50.1B.FI.0S.016.00 SHT UNLOAD X2637'
$ws.Range("E763").Value = '2pcs'
$ws.Range("F763").Value = 'Fishes'
$ws.Range("G763").Value = '15/10/2025'
$ws.Range("H763").Value = 600

# Row 764
$ws.Range("B764").Value = 'SAP'
$ws.Range("C764").Value = 11156538
$ws.Range("D764").Value = 'PTS-1151 TW.50.1B.FI.0S.014.00 SHT LOAD X2637'
$ws.Range("E764").Value = '2pcs'
$ws.Range("F764").Value = 'Fishes'
$ws.Range("G764").Value = '15/10/2025'
$ws.Range("H764").Value = 600

# Row 765
$ws.Range("B765").Value = 'Expense'
$ws.Range("C765").Value = 'Expense'
$ws.Range("D765").Value = '3M Wire Seal Tape (Black)'
$ws.Range("E765").Value = '10roll'
$ws.Range("F765").Value = 'Fishes'
$ws.Range("G765").Value = '16/10/2025'
$ws.Range("H765").Value = 11.58

# Row 766
$ws.Range("B766").Value = 'SAP'
$ws.Range("C766").Value = 11155497
$ws.Range("D766").Value = 'PTS-1138 6K-9090-HS01  LF Vacuum Adaptor'
$ws.Range("E766").Value = '100pcs'
$ws.Range("F766").Value = 'Fishes'
$ws.Range("G766").Value = '18/10/2025'
$ws.Range("H766").Value = 3510

# Row 767
$ws.Range("B767").Value = 'SAP'
$ws.Range("C767").Value = 11155797
$ws.Range("D767").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E767").Value = '100pcs'
$ws.Range("F767").Value = 'Fishes'
$ws.Range("G767").Value = '17/10/2025'
$ws.Range("H767").Value = 2016

# Row 768
$ws.Range("B768").Value = 'SAP'
$ws.Range("C768").Value = 11151246
$ws.Range("D768").Value = 'PTS-1069 TW.50.15.FI.0S.150.00 X2637'
$ws.Range("E768").Value = '3pcs'
$ws.Range("F768").Value = 'Fishes'
$ws.Range("G768").Value = '17/10/2025'
$ws.Range("H768").Value = 555.55999999999995

# Row 769
$ws.Range("B769").Value = 'SAP'
$ws.Range("C769").Value = 11151248
$ws.Range("D769").Value = 'PTS-1071 TW.50.15.FI.0S.151.00 X2637'
$ws.Range("E769").Value = '1pcs'
$ws.Range("F769").Value = 'Fishes'
$ws.Range("G769").Value = '18/10/2025'
$ws.Range("H769").Value = 1666.68

# Row 770
$ws.Range("B770").Value = 'SAP'
$ws.Range("C770").Value = 11155497
$ws.Range("D770").Value = 'PTS-1138 6K-9090-HS01  LF Vacuum Adaptor'
$ws.Range("E770").Value = '97pcs'
$ws.Range("F770").Value = 'Fishes'
$ws.Range("G770").Value = '19/10/2025'
$ws.Range("H770").Value = 3404.7

# Row 771
$ws.Range("B771").Value = 'SAP'
$ws.Range("C771").Value = 11155797
$ws.Range("D771").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E771").Value = '100pcs'
$ws.Range("F771").Value = 'Fishes'
$ws.Range("G771").Value = '19/10/2025'
$ws.Range("H771").Value = 2016

# --- Step 3: fill in column A ("Wk42") for the whole new block last, matching
#     the original authoring order (new "Wk42" shared string is added last). ---
$ws.Range("A722:A771").Value = "Wk42"

# --- Step 4: re-stretch the AutoFilter over the new data range (A1:H771). Excel
#     toggles AutoFilter off if re-applied at the same anchor, so drop + reapply. ---
$ws.AutoFilterMode = $false
$ws.Range("A1:H771").AutoFilter()

# --- Step 5: keep the hidden _FilterDatabase defined name in sync with the filter. ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Weekly Expenditure!_FilterDatabase") {
        $n.RefersTo = "='Weekly Expenditure'!`$A`$1:`$H`$771"
    }
}

# --- Step 6: move the saved selection cursor to match the source workbook. ---
$ws.Range("B724").Select()
